$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the existing hyperlink, currently anchored at C10 (it will need to
# move to C9 once the row above it is deleted).
$ws.Range("C10").Hyperlinks.Delete()

# Delete the old row 3 entirely (Date 2020-12-26 / Postcode 3803 entry).
# This shifts rows 4-28 up to become rows 3-27, shrinking the table from
# A1:E28 to A1:E27.
$ws.Rows("3:3").Delete()

# Correct the date in row 2 to 2020-12-21 (serial 44186), matching the
# corrected "Date" for that entry.
$ws.Range("A2").Value = 44186

# Re-add the hyperlink on the row that shifted from 10 down to 9.
$ws.Hyperlinks.Add($ws.Range("C9"), "https://www.dhhs.vic.gov.au/coronavirus-update-victoria-17-october-2020") | Out-Null

# Move the active selection to A3, mirroring where the edit was made.
$ws.Range("A3").Select() | Out-Null
